$d = $word.ActiveDocument

# Change 1: client step 1 text gains "+ vi" note and a trailing clause.
$d.Content.Find.Execute(
    "+ 1. Генерируется случайный ключ для алгоритма симметричного шифрования (AES, размер ключа на свое усмотрение).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "+ 1. Генерируется случайный ключ для алгоритма симметричного шифрования (AES + vi, размер ключа на свое усмотрение). - file with two Lines + strsplit for the Line for the server;",
    2)

# Change 2: client step 2 gains an extra leading space after the dash.
$d.Content.Find.Execute(
    "- 2. Создается пакет данных - ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "-  2. Создается пакет данных - ",
    2)

# Change 3: client step 5 flips from "-" to "+".
$d.Content.Find.Execute(
    "- 5. Симметричный ключ шифруется ассиметричным шифром (ключи сгенерированы заранее и прописаны в коде программы как константа, алгоритм – RSA, параметры – на свое усмотрение).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "+ 5. Симметричный ключ шифруется ассиметричным шифром (ключи сгенерированы заранее и прописаны в коде программы как константа, алгоритм – RSA, параметры – на свое усмотрение).",
    2)

# Change 4: "Условия" bullet gains "- 8." numbering prefix.
$d.Content.Find.Execute(
    "Выход из приложения по ctrl-c с выводом количества отправленных сообщений в консоль.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- 8. Выход из приложения по ctrl-c с выводом количества отправленных сообщений в консоль.",
    2)

# Change 5: next bullet gains "- 9." numbering prefix.
$d.Content.Find.Execute(
    "Количество одновременно запущенных приложений не ограничено.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- 9. Количество одновременно запущенных приложений не ограничено.",
    2)
